$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new data values (net skip threshold rows)
$ws.Range("B10").Value = 12
$ws.Range("C10").Value = 15
$ws.Range("B11").Value = 12
$ws.Range("C11").Value = 15

# Update the active selection to match the edited location
$ws.Range("C9").Select()
